$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 5.685592333333333
$ws.Range("H2").Value = 17.056777
$ws.Range("I2").Value = 0.1188473284691575
$ws.Range("J2").Value = 0.1188473284691575
$ws.Range("M2").Value = 139.728498
$ws.Range("N2").Value = 419.185494
$ws.Range("O2").Value = 0.9065295391216045
$ws.Range("P2").Value = 0.9065295391216045
$ws.Range("Q2").Value = 794.439276976982
$ws.Range("R2").Value = 7149.953492792838
$ws.Range("S2").Value = 0.1077386139029793
$ws.Range("T2").Value = 0.1077386139029793

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 5.685592333333333
$ws.Range("H3").Value = 17.056777
$ws.Range("I3").Value = 0.1188473284691575
$ws.Range("J3").Value = 0.1188473284691575
$ws.Range("O3").Value = 0.005362677585431591
$ws.Range("P3").Value = 0.005362677585431591
$ws.Range("Q3").Value = 4.699595015689223
$ws.Range("R3").Value = 42.29635514120301
$ws.Range("S3").Value = 0.000637339904469977
$ws.Range("T3").Value = 0.000637339904469977

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 5.685592333333333
$ws.Range("H4").Value = 17.056777
$ws.Range("I4").Value = 0.1188473284691575
$ws.Range("J4").Value = 0.1188473284691575
$ws.Range("O4").Value = 0.0881077832929639
$ws.Range("P4").Value = 0.0881077832929639
$ws.Range("Q4").Value = 77.21346148646278
$ws.Range("R4").Value = 694.9211533781649
$ws.Range("S4").Value = 0.01047137466170823
$ws.Range("T4").Value = 0.01047137466170823

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.622926875404983
$ws.Range("J5").Value = 0.6229268754049829
$ws.Range("M5").Value = 139.728498
$ws.Range("N5").Value = 419.185494
$ws.Range("O5").Value = 0.9065295391216045
$ws.Range("P5").Value = 0.9065295391216045
$ws.Range("Q5").Value = 4163.977288178527
$ws.Range("R5").Value = 37475.79559360674
$ws.Range("S5").Value = 0.5647016132673405
$ws.Range("T5").Value = 0.5647016132673404

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.622926875404983
$ws.Range("J6").Value = 0.6229268754049829
$ws.Range("O6").Value = 0.005362677585431591
$ws.Range("P6").Value = 0.005362677585431591
$ws.Range("S6").Value = 0.00334055599209724
$ws.Range("T6").Value = 0.003340555992097239

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.622926875404983
$ws.Range("J7").Value = 0.6229268754049829
$ws.Range("O7").Value = 0.0881077832929639
$ws.Range("P7").Value = 0.0881077832929639
$ws.Range("S7").Value = 0.05488470614554537
$ws.Range("T7").Value = 0.05488470614554536

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.2582257961258595
$ws.Range("J8").Value = 0.2582257961258594
$ws.Range("M8").Value = 139.728498
$ws.Range("N8").Value = 419.185494
$ws.Range("O8").Value = 0.9065295391216045
$ws.Range("P8").Value = 0.9065295391216045
$ws.Range("Q8").Value = 1726.1196983849
$ws.Range("R8").Value = 15535.0772854641
$ws.Range("S8").Value = 0.2340893119512848
$ws.Range("T8").Value = 0.2340893119512848

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.2582257961258595
$ws.Range("J9").Value = 0.2582257961258594
$ws.Range("O9").Value = 0.005362677585431591
$ws.Range("P9").Value = 0.005362677585431591
$ws.Range("S9").Value = 0.001384781688864374
$ws.Range("T9").Value = 0.001384781688864374

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.2582257961258595
$ws.Range("J10").Value = 0.2582257961258594
$ws.Range("O10").Value = 0.0881077832929639
$ws.Range("P10").Value = 0.0881077832929639
$ws.Range("S10").Value = 0.0227517024857103
$ws.Range("T10").Value = 0.0227517024857103

Write-Output "Applied TPM updates to Adam9-Itga6 sheet"
